{"js": "// Update the heading run to expand \"release 1\" -> \"release 4 ... on the\n// external emulator\" (splitting the single run into several, matching the\n// author's edit), and relocate the stray \"_GoBack\" bookmark from the end of\n// the document to the end of that same heading line.\n\nconst body = context.document.body;\n\n// 1) Remove the \"_GoBack\" bookmark from its old spot (right after\n//    \"To rotate the shape, use the 'w' key\"). Must happen BEFORE the new\n//    bookmark is inserted below, otherwise the document would contain two\n//    bookmarks with the same name and a lookup-by-name would not\n//    deterministically hit the trailing one anymore.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Find the heading run that needs to be rewritten.\nconst results = body.search(\n  \"Instructions for running the full release 1 code for Tetris and Pong:\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const rPr = '<w:rPr><w:u w:val=\"single\"/><w:lang w:val=\"en-CA\"/></w:rPr>';\n  const flatOpc = (inner) =>\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + inner + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n\n  const newContent =\n    '<w:p>' +\n      '<w:r>' + rPr + '<w:t xml:space=\"preserve\">Instructions for running the full release </w:t></w:r>' +\n      '<w:r>' + rPr + '<w:t>4</w:t></w:r>' +\n      '<w:r>' + rPr + '<w:t xml:space=\"preserve\"> code for Tetris and Pong</w:t></w:r>' +\n      '<w:r>' + rPr + '<w:t xml:space=\"preserve\"> on the external emulator</w:t></w:r>' +\n      '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n      '<w:r>' + rPr + '<w:t>:</w:t></w:r>' +\n    '</w:p>';\n\n  results.items[0].insertOoxml(flatOpc(newContent), \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the heading run to expand \"release 1\" -> \"release 4 ... on the\n# external emulator\" (splitting the single run into several, matching the\n# author's edit), and relocate the stray \"_GoBack\" bookmark from the end of\n# the document to the end of that same heading line.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the \"_GoBack\" bookmark from its old spot (right after\n#    \"To rotate the shape, use the 'w' key\"). This must happen BEFORE the\n#    new bookmark is inserted below, otherwise the document would contain\n#    two bookmarks named \"_GoBack\" and a name lookup could land on the\n#    stale one instead of the new one.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Find the heading run that needs to be rewritten.\n$range = $d.Content\n$found = $range.Find.Execute(\"Instructions for running the full release 1 code for Tetris and Pong:\")\n\nif ($found) {\n  $rPr = '<w:rPr><w:u w:val=\"single\"/><w:lang w:val=\"en-CA\"/></w:rPr>'\n  $newContent = '<w:p>' +\n    '<w:r>' + $rPr + '<w:t xml:space=\"preserve\">Instructions for running the full release </w:t></w:r>' +\n    '<w:r>' + $rPr + '<w:t>4</w:t></w:r>' +\n    '<w:r>' + $rPr + '<w:t xml:space=\"preserve\"> code for Tetris and Pong</w:t></w:r>' +\n    '<w:r>' + $rPr + '<w:t xml:space=\"preserve\"> on the external emulator</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r>' + $rPr + '<w:t>:</w:t></w:r>' +\n    '</w:p>'\n\n  $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + $newContent + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n  $range.InsertXML($xml)\n}\n"}
